$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename county "Aston" -> "Delaware County" for the three rows that used it (5, 6, 7)
$ws.Range("D5").Value = "Delaware County"
$ws.Range("D6").Value = "Delaware County"
$ws.Range("D7").Value = "Delaware County"

# Rename county "Camden" -> "Camden County" for row 8
$ws.Range("D8").Value = "Camden County"

# Update median_household_income (column R) values
$ws.Range("R2").Value = 85000
$ws.Range("R4").Value = 85000
$ws.Range("R6").Value = 80000
$ws.Range("R7").Value = 80000
$ws.Range("R8").Value = 67000

# Move / collapse the active selection to D8
$ws.Range("D8").Select()
